$wb = $excel.ActiveWorkbook

# --- Sheets ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (in zh-cn / de-de tables) whose "Priority" changes from blank to "ht"
# (these are the rows for: 48a808f8, 69ce540c, 7d3a2894, bda6caee, c25b16a1, dc0e7a20)
$rows = @(7, 9, 10, 11, 12, 14)

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}

# New handoff timestamps generated for this handoff report
$newDeDate = "2016-08-19 04:18:33"   # de-de handoff datetime / overview's latest HO xliff generate date
$newZhDate = "2016-08-19 04:18:28"   # zh-cn handoff datetime

foreach ($r in $rows) {
    $overview.Range("G$r").Value = $newDeDate
    $dede.Range("H$r").Value     = $newDeDate
    $zhcn.Range("H$r").Value     = $newZhDate
}
